$d = $word.ActiveDocument

# Replace "final year" with "final semester" in the cover-letter paragraph.
$rng = $d.Content
$rng.Find.Execute("final year", $true, $false, $false, $false, $false, $true, 1, $false, "final semester", 2)

# Move the "_GoBack" bookmark so it sits right after the word "semester"
# (matches Word's behaviour of re-anchoring _GoBack at the last edit point).
$pos = $rng.End
$bmRange = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $bmRange)
